$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.490.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.58%  '
$ws.Range("D3").Value = "'2.476.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.47%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = "'573.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.91%  '
$ws.Range("D6").Value = "'149.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.12%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = "'0.542"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.33%  '
$ws.Range("E9").Value = '  +4.58%  '
$ws.Range("E10").Value = '  +0.49%  '
$ws.Range("D11").Value = "'0.364"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.35%  '
$ws.Range("E12").Value = '  +2.58%  '
$ws.Range("D13").Value = "'27.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.70%  '
$ws.Range("E14").Value = '  +7.08%  '
$ws.Range("D15").Value = "'2.918.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.53%  '
$ws.Range("D16").Value = "'63.334.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.46%  '
$ws.Range("D17").Value = "'2.484.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.88%  '
$ws.Range("D18").Value = "'11.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.66%  '
$ws.Range("D19").Value = "'7.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.71%  '
$ws.Range("E20").Value = '  +3.34%  '
$ws.Range("D21").Value = "'328.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.58%  '
$ws.Range("D22").Value = "'0.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("E23").Value = '  +11.11%  '
$ws.Range("D24").Value = "'67.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.67%  '
$ws.Range("D25").Value = "'636.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +14.87%  '
$ws.Range("D26").Value = "'8.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.19%  '
$ws.Range("E27").Value = '  +12.98%  '
$ws.Range("D28").Value = "'2.603.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.63%  '
$ws.Range("E29").Value = '  +9.82%  '
$ws.Range("D30").Value = "'8.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.17%  '
$ws.Range("D31").Value = "'0.994"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("E32").Value = '  -1.68%  '
$ws.Range("D33").Value = "'1.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.92%  '
$ws.Range("E34").Value = '  +10.86%  '
$ws.Range("E35").Value = '  +4.71%  '
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("E37").Value = '  +2.31%  '
$ws.Range("D38").Value = "'5.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.58%  '
$ws.Range("D39").Value = "'18.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.32%  '
$ws.Range("E40").Value = '  +2.67%  '
$ws.Range("D41").Value = "'146.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.56%  '
$ws.Range("D42").Value = "'2.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +18.67%  '
$ws.Range("E43").Value = '  +0.85%  '
$ws.Range("D44").Value = "'151.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.19%  '
$ws.Range("D45").Value = "'3.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.59%  '
$ws.Range("D46").Value = "'0.0552"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.16%  '
$ws.Range("D47").Value = "'21.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.15%  '
$ws.Range("E48").Value = '  +3.45%  '
$ws.Range("E49").Value = '  +6.49%  '
$ws.Range("D50").Value = "'0.0927"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.89%  '
$ws.Range("D51").Value = "'0.744"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.83%  '
